# Updates crypto price/volume data in Sheet1 to reflect the latest scrape
# (GitHub Actions scheduled refresh). Also corrects the ordering of the
# "ordi" / "Algorand" rows (47 and 48), whose coin name, link, price, and
# volume values were swapped in the source feed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '46.775.78'
$ws.Range("E2").Value = '  -0.42%  '
$ws.Range("D3").Value = '2.261.89'
$ws.Range("E3").Value = '  -3.74%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = '298.08'
$ws.Range("E5").Value = '  -2.52%  '
$ws.Range("D6").Value = '97.62'
$ws.Range("E6").Value = '  -0.50%  '
$ws.Range("D7").Value = '0.573'
$ws.Range("E7").Value = '  -0.66%  '
$ws.Range("E8").Value = '  +0.09%  '
$ws.Range("E9").Value = '  -7.04%  '
$ws.Range("D10").Value = '34.55'
$ws.Range("E10").Value = '  -3.48%  '
$ws.Range("D11").Value = '0.0788'
$ws.Range("E11").Value = '  -2.20%  '
$ws.Range("D12").Value = '6.97'
$ws.Range("E12").Value = '  -6.36%  '
$ws.Range("E13").Value = '  -1.78%  '
$ws.Range("D14").Value = '2.605.91'
$ws.Range("E14").Value = '  -3.72%  '
$ws.Range("D15").Value = '2.260.08'
$ws.Range("E15").Value = '  -3.64%  '
$ws.Range("D16").Value = '13.54'
$ws.Range("E16").Value = '  -4.76%  '
$ws.Range("D17").Value = '46.713.59'
$ws.Range("E17").Value = '  -0.19%  '
$ws.Range("E18").Value = '  -5.06%  '
$ws.Range("D19").Value = '0.0₃0964'
$ws.Range("E19").Value = '  +1.53%  '
$ws.Range("D20").Value = '12.31'
$ws.Range("E20").Value = '  -10.57%  '
$ws.Range("E21").Value = '  -6.93%  '
$ws.Range("D22").Value = '65.66'
$ws.Range("E22").Value = '  -1.77%  '
$ws.Range("D23").Value = '243.54'
$ws.Range("E23").Value = '  -0.62%  '
$ws.Range("D24").Value = '2.76'
$ws.Range("E24").Value = '  -7.52%  '
$ws.Range("E25").Value = '  +0.21%  '
$ws.Range("E26").Value = '  -8.02%  '
$ws.Range("D27").Value = '40.91'
$ws.Range("E27").Value = '  -2.03%  '
$ws.Range("E28").Value = '  -3.86%  '
$ws.Range("D29").Value = '9.47'
$ws.Range("E29").Value = '  -4.28%  '
$ws.Range("D30").Value = '19.99'
$ws.Range("E30").Value = '  -0.88%  '
$ws.Range("D31").Value = '2.81'
$ws.Range("E31").Value = '  +7.15%  '
$ws.Range("D32").Value = '3.29'
$ws.Range("E32").Value = '  +3.83%  '
$ws.Range("D33").Value = '144.05'
$ws.Range("E33").Value = '  -5.51%  '
$ws.Range("D34").Value = '5.26'
$ws.Range("E34").Value = '  -8.91%  '
$ws.Range("D35").Value = '0.0761'
$ws.Range("E35").Value = '  -6.97%  '
$ws.Range("E36").Value = '  +1.03%  '
$ws.Range("D38").Value = '15.19'
$ws.Range("E38").Value = '  +10.00%  '
$ws.Range("D39").Value = '1.65'
$ws.Range("E39").Value = '  -10.45%  '
$ws.Range("D40").Value = '3.78'
$ws.Range("E40").Value = '  -6.40%  '
$ws.Range("E41").Value = '  -7.47%  '
$ws.Range("E42").Value = '  -10.99%  '
$ws.Range("D43").Value = '0.998'
$ws.Range("E43").Value = '  -0.04%  '
$ws.Range("D44").Value = '92.13'
$ws.Range("E44").Value = '  +13.65%  '
$ws.Range("D45").Value = '1.777.94'
$ws.Range("E45").Value = '  -4.63%  '
$ws.Range("D46").Value = '1.85'
$ws.Range("E46").Value = '  -7.14%  '
$ws.Range("B47").Value = 'Algorand'
$ws.Range("C47").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D47").Value = '0.182'
$ws.Range("E47").Value = '  -7.75%  '
$ws.Range("B48").Value = 'ordi'
$ws.Range("C48").Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range("D48").Value = '69.54'
$ws.Range("E48").Value = '  -6.46%  '
$ws.Range("D49").Value = '4.76'
$ws.Range("E49").Value = '  -3.77%  '
$ws.Range("D50").Value = '2.484.23'
$ws.Range("E50").Value = '  -3.79%  '
$ws.Range("E51").Value = '  -5.29%  '
